$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: assigning a bare "" to Range.Value clears/removes the cell entirely.
# To leave the cell in place as an explicit empty text value (matching the
# target workbook, which keeps these cells present with blank text), a
# lone leading apostrophe ("'") is used - Excel's standard "force text"
# input marker. It is stripped on input, leaving a real (non-null) empty
# string value of type Text rather than an absent cell.

# Row 2
$ws.Range("M2").Value = "'"
$ws.Range("N2").Value = "'"
$ws.Range("O2").Value = "'"
$ws.Range("Z2").Value = "Suppress"
$ws.Range("AA2").Value = "'"
$ws.Range("AB2").Value = "1. Video Version number is not available in CAS."
$ws.Range("AC2").Value = "Dummy, Title"

# Row 3
$ws.Range("M3").Value = "'"
$ws.Range("N3").Value = "'"
$ws.Range("O3").Value = "'"
$ws.Range("Z3").Value = "'"
$ws.Range("AA3").NumberFormat = "@"
$ws.Range("AA3").Value = "123"
$ws.Range("AB3").Value = "1. Video Version number is not available in CAS."
$ws.Range("AC3").Value = "Dummy, Title1"
